# Generate Report for Handback
#
# A new handback entry (5b8f047b-a03d-4a5f-b517-c61587b9de23.md) is inserted
# as the *second* data row (row 3) on every sheet, pushing the existing
# 5d7baab8-... entry down to a brand-new last row (row 4). Each of the three
# tables (Overview, zh-cn, de-de) grows by one row.
#
# NB: every textual value is written with a leading "'" so the engine keeps
# it typed as a shared string (t="s") instead of auto-coercing look-alike
# values ("True"/"False", empty string, date-looking text) into booleans /
# numbers - that's how the source workbook stores them too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data describing the new / shifted entries
# ---------------------------------------------------------------------------
$newGuid = "5b8f047b-a03d-4a5f-b517-c61587b9de23"
$newMd = "$newGuid.md"
$newDisplayBackslash = "e2e\$newGuid.md"

$oldGuid = "5d7baab8-2b3d-445c-a719-9c3245fc4841"
$oldMd = "$oldGuid.md"
$oldDisplayBackslash = "e2e\$oldGuid.md"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-Text($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

function Set-DateText($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $text
    $r.NumberFormat = $dateFmt
}

function Remove-HyperlinkAt($ws, $addr) {
    $target = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $target = $hl
        }
    }
    if ($target -ne $null) {
        $target.Delete()
    }
}

# ---------------------------------------------------------------------------
# 1) "Overview" sheet - columns A..G, table3.xml (displayName "Overview")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# Snapshot the current (old) row 3 -- this is the 5d7baab8 entry which needs
# to move down to row 4.
$ov_a3 = $ws.Range("A3").Value()
$ov_c3 = $ws.Range("C3").Value()
$ov_e3 = $ws.Range("E3").Value()
$ov_f3 = $ws.Range("F3").Value()
$ov_g3 = $ws.Range("G3").Value()

# Grow the table by one row (new row becomes row 4).
$lo.ListRows.Add() | Out-Null

# Move the old row 3 content down into the new row 4.
Set-Text $ws "A4" $ov_a3
Set-Text $ws "C4" $ov_c3
Set-Text $ws "E4" $ov_e3
Set-Text $ws "F4" $ov_f3
Set-DateText $ws "G4" $ov_g3

$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3dffdf356cb27dddc574f78e7dd5a20d5d1630fa/e2e/$oldMd", "", "", $oldDisplayBackslash) | Out-Null

# Replace row 3 with the new 5b8f047b entry.
Remove-HyperlinkAt $ws '$B$3'

Set-Text $ws "A3" $newMd
Set-Text $ws "C3" ".md"
Set-Text $ws "E3" "Handed back: in sync with en-US"
Set-Text $ws "F3" "Handed back: in sync with en-US"
Set-DateText $ws "G3" "2016-09-02 06:52:14"

$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fec46ebfc87d6af1ca4b76d79e8bda98b9b893c6/e2e/$newMd", "", "", $newDisplayBackslash) | Out-Null

# ---------------------------------------------------------------------------
# 2) "zh-cn" sheet - columns A..P, table1.xml
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

$zh_b3 = $ws.Range("B3").Value()
$zh_c3 = $ws.Range("C3").Value()
$zh_d3 = $ws.Range("D3").Value()
$zh_e3 = $ws.Range("E3").Value()
$zh_f3 = $ws.Range("F3").Value()
$zh_g3 = $ws.Range("G3").Value()
$zh_h3 = $ws.Range("H3").Value()
$zh_j3 = $ws.Range("J3").Value()
$zh_k3 = $ws.Range("K3").Value()
$zh_l3 = $ws.Range("L3").Value()
$zh_m3 = $ws.Range("M3").Value()
$zh_n3 = $ws.Range("N3").Value()
$zh_o3 = $ws.Range("O3").Value()
$zh_p3 = $ws.Range("P3").Value()

$lo.ListRows.Add() | Out-Null

Set-Text $ws "A4" $oldMd
Set-Text $ws "B4" $zh_b3
Set-Text $ws "C4" $zh_c3
Set-Text $ws "D4" $zh_d3
Set-Text $ws "E4" $zh_e3
Set-Text $ws "F4" $zh_f3
Set-Text $ws "G4" $zh_g3
Set-DateText $ws "H4" $zh_h3
Set-Text $ws "I4" $oldMd
Set-Text $ws "J4" $zh_j3
Set-DateText $ws "K4" $zh_k3
Set-Text $ws "L4" $zh_l3
Set-Text $ws "M4" $zh_m3
Set-Text $ws "N4" $zh_n3
Set-Text $ws "O4" $zh_o3
Set-Text $ws "P4" $zh_p3

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3dffdf356cb27dddc574f78e7dd5a20d5d1630fa/e2e/$oldMd", "", "", $oldMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/881de67ccdcf9b438fa4c9ab05d267efb450be1b/e2e/$oldMd", "", "", $oldMd) | Out-Null

Remove-HyperlinkAt $ws '$A$3'
Remove-HyperlinkAt $ws '$I$3'

Set-Text $ws "A3" $newMd
Set-Text $ws "B3" ".md"
Set-Text $ws "C3" "Handed back: in sync with en-US"
Set-Text $ws "D3" "e2e"
Set-Text $ws "E3" "ht"
Set-Text $ws "F3" "True"
Set-Text $ws "G3" "$newGuid.5959d13bab1e9e57f919f4850e88d6bc590264f9.zh-cn.xlf"
Set-DateText $ws "H3" "2016-09-02 06:52:08"
Set-Text $ws "I3" $newMd
Set-Text $ws "J3" "$newGuid.5959d13bab1e9e57f919f4850e88d6bc590264f9.zh-cn.xlf"
Set-DateText $ws "K3" "2016-09-02 06:52:36"
Set-Text $ws "L3" ""
Set-Text $ws "M3" "True"
Set-Text $ws "N3" ""
Set-Text $ws "O3" "False"
Set-Text $ws "P3" ""

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fec46ebfc87d6af1ca4b76d79e8bda98b9b893c6/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/43a472a99dd04823f20f4f8cff831eac4f928565/e2e/$newMd", "", "", $newMd) | Out-Null

# ---------------------------------------------------------------------------
# 3) "de-de" sheet - columns A..P, table2.xml
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

$de_b3 = $ws.Range("B3").Value()
$de_c3 = $ws.Range("C3").Value()
$de_d3 = $ws.Range("D3").Value()
$de_e3 = $ws.Range("E3").Value()
$de_f3 = $ws.Range("F3").Value()
$de_g3 = $ws.Range("G3").Value()
$de_h3 = $ws.Range("H3").Value()
$de_j3 = $ws.Range("J3").Value()
$de_k3 = $ws.Range("K3").Value()
$de_l3 = $ws.Range("L3").Value()
$de_m3 = $ws.Range("M3").Value()
$de_n3 = $ws.Range("N3").Value()
$de_o3 = $ws.Range("O3").Value()
$de_p3 = $ws.Range("P3").Value()

$lo.ListRows.Add() | Out-Null

Set-Text $ws "A4" $oldMd
Set-Text $ws "B4" $de_b3
Set-Text $ws "C4" $de_c3
Set-Text $ws "D4" $de_d3
Set-Text $ws "E4" $de_e3
Set-Text $ws "F4" $de_f3
Set-Text $ws "G4" $de_g3
Set-DateText $ws "H4" $de_h3
Set-Text $ws "I4" $oldMd
Set-Text $ws "J4" $de_j3
Set-DateText $ws "K4" $de_k3
Set-Text $ws "L4" $de_l3
Set-Text $ws "M4" $de_m3
Set-Text $ws "N4" $de_n3
Set-Text $ws "O4" $de_o3
Set-Text $ws "P4" $de_p3

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3dffdf356cb27dddc574f78e7dd5a20d5d1630fa/e2e/$oldMd", "", "", $oldMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f727f3396834ee3237a5c1a381f9cd5c3616fce5/e2e/$oldMd", "", "", $oldMd) | Out-Null

Remove-HyperlinkAt $ws '$A$3'
Remove-HyperlinkAt $ws '$I$3'

Set-Text $ws "A3" $newMd
Set-Text $ws "B3" ".md"
Set-Text $ws "C3" "Handed back: in sync with en-US"
Set-Text $ws "D3" "e2e"
Set-Text $ws "E3" "ht"
Set-Text $ws "F3" "True"
Set-Text $ws "G3" "$newGuid.5959d13bab1e9e57f919f4850e88d6bc590264f9.de-de.xlf"
Set-DateText $ws "H3" "2016-09-02 06:52:14"
Set-Text $ws "I3" $newMd
Set-Text $ws "J3" "$newGuid.5959d13bab1e9e57f919f4850e88d6bc590264f9.de-de.xlf"
Set-DateText $ws "K3" "2016-09-02 06:52:43"
Set-Text $ws "L3" ""
Set-Text $ws "M3" "True"
Set-Text $ws "N3" ""
Set-Text $ws "O3" "False"
Set-Text $ws "P3" ""

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fec46ebfc87d6af1ca4b76d79e8bda98b9b893c6/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7898e9fd7da8643c60ea3b299d21149869959cf3/e2e/$newMd", "", "", $newMd) | Out-Null
